$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "Miss Dina Nasr, Administrator"
$newText = "Administrator, Miss Dina Nasr"

# The credentials/role text in column G ("Miss Dina Nasr, Administrator")
# is reordered to "Administrator, Miss Dina Nasr" across every row of the
# session analysis table. Scan the full used range of column G and swap
# any cell whose text matches exactly.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
